# Scheduled market-data refresh: update computed price/profit columns
# (H..N) on the per-profession "Leve" sheets. Values below are the
# freshly recalculated currentAveragePrice / LevePrice / LeveProfit
# figures pulled for each affected Leve row; unrelated columns (A..G)
# are left untouched.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 155.14285
$ws.Range("I4").Value = 155.14285
$ws.Range("K4").Value = 155.14285
$ws.Range("M4").Value = -41.14285000000001

$ws.Range("H7").Value = 6250.25
$ws.Range("I7").Value = 6000.3335
$ws.Range("J7").Value = 6400.2
$ws.Range("K7").Value = 6000.3335
$ws.Range("L7").Value = 6400.2
$ws.Range("M7").Value = -5888.3335
$ws.Range("N7").Value = -6624.2

$ws.Range("H14").Value = 6250.25
$ws.Range("I14").Value = 6000.3335
$ws.Range("J14").Value = 6400.2
$ws.Range("K14").Value = 6000.3335
$ws.Range("L14").Value = 6400.2
$ws.Range("M14").Value = -5809.3335
$ws.Range("N14").Value = -6782.2

$ws.Range("H17").Value = 2360.818
$ws.Range("J17").Value = 2360.818
$ws.Range("L17").Value = 7082.454000000001
$ws.Range("N17").Value = -7418.454000000001

$ws.Range("H70").Value = 2770.375
$ws.Range("J70").Value = 3657.4
$ws.Range("L70").Value = 10972.2
$ws.Range("N70").Value = -11512.2

$ws.Range("H73").Value = 2770.375
$ws.Range("J73").Value = 3657.4
$ws.Range("L73").Value = 10972.2
$ws.Range("N73").Value = -12844.2

$ws.Range("H86").Value = 6732.6665
$ws.Range("I86").Value = 6732.6665
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 6732.6665
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -5609.6665
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 6732.6665
$ws.Range("I89").Value = 6732.6665
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 33663.3325
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -28047.3325
$ws.Range("N89").ClearContents()

$ws.Range("H93").Value = 28532.334
$ws.Range("J93").Value = 28532.334
$ws.Range("L93").Value = 28532.334
$ws.Range("N93").Value = -33524.334

$ws.Range("H100").Value = 700.875
$ws.Range("J100").Value = 287
$ws.Range("L100").Value = 287
$ws.Range("N100").Value = -1369

$ws.Range("H103").Value = 536.4
$ws.Range("J103").Value = 536.4
$ws.Range("L103").Value = 1609.2
$ws.Range("N103").Value = -2781.2

$ws.Range("H107").Value = 2784.6667
$ws.Range("I107").Value = 722.4286
$ws.Range("K107").Value = 722.4286
$ws.Range("M107").Value = 1197.5714

$ws.Range("H112").Value = 1842.125
$ws.Range("J112").Value = 2014.5834
$ws.Range("L112").Value = 6043.7502
$ws.Range("N112").Value = -8259.7502

$ws.Range("H132").Value = 2224.524
$ws.Range("I132").Value = 2285.8
$ws.Range("K132").Value = 6857.400000000001
$ws.Range("M132").Value = -4327.400000000001

$ws.Range("H137").Value = 3652.6924
$ws.Range("I137").Value = 2164.6667
$ws.Range("J137").Value = 4099.1
$ws.Range("K137").Value = 6494.000100000001
$ws.Range("L137").Value = 12297.3
$ws.Range("M137").Value = -3944.000100000001
$ws.Range("N137").Value = -17397.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1792.8
$ws.Range("I45").Value = 1632.4445
$ws.Range("K45").Value = 1632.4445
$ws.Range("M45").Value = -1255.4445

$ws.Range("H95").Value = 12999.25
$ws.Range("J95").Value = 11998.5
$ws.Range("L95").Value = 11998.5
$ws.Range("N95").Value = -17490.5

$ws.Range("H102").Value = 2369
$ws.Range("I102").Value = 1597.1666
$ws.Range("J102").Value = 7000
$ws.Range("K102").Value = 1597.1666
$ws.Range("L102").Value = 7000
$ws.Range("M102").Value = 24.83339999999998
$ws.Range("N102").Value = -10244

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6999
$ws.Range("I20").Value = 6999.5
$ws.Range("J20").Value = 6998
$ws.Range("K20").Value = 6999.5
$ws.Range("L20").Value = 6998
$ws.Range("M20").Value = -6752.5
$ws.Range("N20").Value = -7492

$ws.Range("H99").Value = 2732
$ws.Range("I99").Value = 2732
$ws.Range("K99").Value = 2732
$ws.Range("M99").Value = -1234

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2104.8572
$ws.Range("I105").Value = 1250
$ws.Range("K105").Value = 1250
$ws.Range("M105").Value = 497

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 85.916664
$ws.Range("I12").Value = 36.666668
$ws.Range("J12").Value = 102.333336
$ws.Range("K12").Value = 110.000004
$ws.Range("L12").Value = 307.000008
$ws.Range("M12").Value = 62.999996
$ws.Range("N12").Value = -653.000008

$ws.Range("H37").Value = 100000
$ws.Range("J37").Value = 100000
$ws.Range("L37").Value = 300000
$ws.Range("N37").Value = -300224

$ws.Range("H68").Value = 1490.3334
$ws.Range("I68").Value = 1623
$ws.Range("K68").Value = 4869
$ws.Range("M68").Value = -4058

$ws.Range("H71").Value = 1490.3334
$ws.Range("I71").Value = 1623
$ws.Range("K71").Value = 14607
$ws.Range("M71").Value = -10551

$ws.Range("H86").Value = 685.1429000000001
$ws.Range("I86").Value = 558.2
$ws.Range("K86").Value = 1674.6
$ws.Range("M86").Value = -488.6000000000001

$ws.Range("H89").Value = 685.1429000000001
$ws.Range("I89").Value = 558.2
$ws.Range("K89").Value = 5023.8
$ws.Range("M89").Value = 904.1999999999998

$ws.Range("H107").Value = 1166.6666
$ws.Range("I107").Value = 1142.8572
$ws.Range("J107").Value = 1250
$ws.Range("K107").Value = 3428.5716
$ws.Range("L107").Value = 3750
$ws.Range("M107").Value = -1508.5716
$ws.Range("N107").Value = -7590

$ws.Range("H112").Value = 1355.6666
$ws.Range("I112").Value = 18.5
$ws.Range("K112").Value = 55.5
$ws.Range("M112").Value = 1052.5

$ws.Range("H113").Value = 1932.6666
$ws.Range("J113").Value = 2079
$ws.Range("L113").Value = 6237
$ws.Range("N113").Value = -10577

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 35802.11
$ws.Range("J20").Value = 35802.11
$ws.Range("L20").Value = 35802.11
$ws.Range("N20").Value = -36292.11

$ws.Range("H80").Value = 4497.5
$ws.Range("J80").Value = 4497.5
$ws.Range("L80").Value = 4497.5
$ws.Range("N80").Value = -6493.5

$ws.Range("H83").Value = 4497.5
$ws.Range("J83").Value = 4497.5
$ws.Range("L83").Value = 22487.5
$ws.Range("N83").Value = -32471.5

$ws.Range("H97").Value = 537.53845
$ws.Range("I97").Value = 353.45456
$ws.Range("J97").Value = 1550
$ws.Range("K97").Value = 353.45456
$ws.Range("L97").Value = 1550
$ws.Range("M97").Value = 142.54544
$ws.Range("N97").Value = -2542

$ws.Range("H98").Value = 7396.25
$ws.Range("J98").Value = 7396.25
$ws.Range("L98").Value = 7396.25
$ws.Range("N98").Value = -13386.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1437.8
$ws.Range("I82").Value = 1347
$ws.Range("J82").Value = 1498.3334
$ws.Range("K82").Value = 1347
$ws.Range("L82").Value = 1498.3334
$ws.Range("M82").Value = -986
$ws.Range("N82").Value = -2220.3334

$ws.Range("H85").Value = 1437.8
$ws.Range("I85").Value = 1347
$ws.Range("J85").Value = 1498.3334
$ws.Range("K85").Value = 1347
$ws.Range("L85").Value = 1498.3334
$ws.Range("M85").Value = -99
$ws.Range("N85").Value = -3994.3334

$ws.Range("H93").Value = 1331
$ws.Range("I93").Value = 996.5
$ws.Range("K93").Value = 996.5
$ws.Range("M93").Value = 251.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 500
$ws.Range("I7").Value = 500
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 500
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -387
$ws.Range("N7").ClearContents()
